# Update the C6 capacitor BOM row (row 12) to the new part:
#  - MFG part number (Comment/A): CL21A226MAQNNNE -> C3216X5R1V226MTJ00E
#  - Value (B): 22uF 25V 0805 -> 22UF 35V 1206
#  - Footprint (D): 0805_C -> 1206_C
$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("A12").Value = "C3216X5R1V226MTJ00E "
$ws.Range("B12").Value = "22UF ±20% 35V X5R 1206 MULTILAYER CERAMIC CAPACITORS MLCC "
$ws.Range("D12").Value = "1206_C"

# Append a new BOM line at row 48, with only the Designator column (B) filled
# in, matching the existing "D5" designator text used elsewhere (C5).
$ws.Range("B48").Value = "D5"

# Match the active-cell selection left behind by the edit.
$ws.Range("E1").Select()
